$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-home existing formatting BEFORE we overwrite values, so cells keep
#    (or pick up) the correct look without minting brand-new style records.
# ---------------------------------------------------------------------------

# A5 needs to become the "label, centered" look that A6 already carries.
$ws.Cells.Item(6, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null

# A6 needs to become a plain label like A2/A3/A4.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(6, 1).PasteSpecial(-4122) | Out-Null

# A8 needs to become the look A9 currently carries (grab it before A9 changes).
$ws.Cells.Item(9, 1).Copy() | Out-Null
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null

# B5 currently carries the old hyperlink-style formatting; make it a plain
# value cell like its neighbours (B2/B6/B7/...).
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Cells.Item(5, 2).PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2) Write the new webDriver / config rows.
# ---------------------------------------------------------------------------

$ws.Cells.Item(5, 1).Value = "IEServerLocation"
$ws.Cells.Item(5, 2).Value = "Resources\servers\IEDriverServer.exe"
$ws.Cells.Item(5, 3).Value = ""

$ws.Cells.Item(6, 1).Value = "ChromeServerLocation"
$ws.Cells.Item(6, 2).Value = "Resources\servers\chromedriver.exe"
$ws.Cells.Item(6, 3).Value = ""

$ws.Cells.Item(7, 1).Value = "PhantomJSLocation"
$ws.Cells.Item(7, 2).Value = "Resources\servers\phantomjs.exe"
$ws.Cells.Item(7, 3).Value = ""

$ws.Cells.Item(8, 1).Value = "ConfigDevUrl"
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(8, 3).Value = ""

# A9: reset to the plain "Normal" look (left/top aligned), drop B9/C9 entirely.
$a9 = $ws.Cells.Item(9, 1)
$a9.Style = "Normal"
$a9.HorizontalAlignment = -4131
$a9.VerticalAlignment = -4160
$a9.Value = "ConfigDevUsername"
$ws.Range("B9:C9").Clear() | Out-Null

# A10: brand new row with the same plain look as A9.
$a10 = $ws.Cells.Item(10, 1)
$a10.Style = "Normal"
$a10.HorizontalAlignment = -4131
$a10.VerticalAlignment = -4160
$a10.Value = "ConfigDevPassword"

# ---------------------------------------------------------------------------
# 3) Sheet-level cosmetics: column width, selection, dimension (dimension is
#    recalculated automatically by the engine from used cells).
# ---------------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 34.6
$ws.Range("A5").Select() | Out-Null
